# Add the second login test case ("login with valid cred") as a new row
# beneath the existing "Launch application" scenario.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test-case row (row 3): slno, Test Scenario, Precondition, Test Steps,
# Expected Result, Test Result.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "login with valid cred"
$ws.Range("C3").Value = "user should be in login scree"
$ws.Range("D3").Value = "Enter valid UserName and Password then click on login button"
$ws.Range("E3").Value = "User should be able to enter user name and password"
$ws.Range("F3").Value = "Fail"

# Widen the "Test Scenario" column to fit the new text.
$ws.Columns(2).ColumnWidth = 17

# Leave the cursor where the author left it after typing the row.
$ws.Range("G3").Select()
